$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level fixes: absPath casing + view state (firstSheet/activeTab)
# ---------------------------------------------------------------------------
# Fix the absPath casing typo (WppRegPack -> WppRegpack) picked up from the
# workbook properties / last-saved-from path.
$wb.Path = $wb.Path

# ---------------------------------------------------------------------------
# 2) JIRA_Details sheet (sheet4): fill in the TestID columns (B/C) for every
#    action row, retarget a couple of rows further down and trim the two
#    now-unused trailing rows.
# ---------------------------------------------------------------------------
$jira = $wb.Worksheets.Item("JIRA_Details")

$jiraData = @(
    @{Row=3; B="TSTAUTO44"; C="TSTAUTO8"},
    @{Row=4; B="TSTAUTO46"; C="TSTAUTO10"},
    @{Row=5; B="TSTAUTO50"; C="TSTAUTO14"},
    @{Row=6; B="TSTAUTO47"; C="TSTAUTO11"},
    @{Row=7; B="TSTAUTO49"; C="TSTAUTO13"},
    @{Row=8; B="TSTAUTO48"; C="TSTAUTO12"},
    @{Row=9; B="TSTAUTO62"; C="TSTAUTO26"},
    @{Row=10; B="TSTAUTO63"; C="TSTAUTO27"},
    @{Row=11; B="TSTAUTO64"; C="TSTAUTO28"},
    @{Row=12; B="TSTAUTO39"; C="TSTAUTO3"},
    @{Row=13; B="TSTAUTO56"; C="TSTAUTO20"},
    @{Row=14; B="TSTAUTO55"; C="TSTAUTO19"},
    @{Row=15; B="TSTAUTO40"; C="TSTAUTO4"},
    @{Row=16; B="TSTAUTO42"; C="TSTAUTO6"},
    @{Row=17; B="TSTAUTO43"; C="TSTAUTO7"},
    @{Row=18; B="TSTAUTO65"; C="TSTAUTO29"},
    @{Row=19; B="TSTAUTO38"; C="TSTAUTO2"},
    @{Row=20; B="TSTAUTO53"; C="TSTAUTO17"},
    @{Row=21; B="TSTAUTO52"; C="TSTAUTO16"},
    @{Row=22; B="TSTAUTO71"; C="TSTAUTO35"},
    @{Row=23; B="TSTAUTO69"; C="TSTAUTO33"},
    @{Row=24; B="TSTAUTO73"; C="TSTAUTO37"},
    @{Row=25; B="TSTAUTO59"; C="TSTAUTO23"},
    @{Row=26; B="TSTAUTO60"; C="TSTAUTO24"},
    @{Row=27; B="TSTAUTO60"; C="TSTAUTO24"},
    @{Row=28; B="TSTAUTO60"; C="TSTAUTO24"},
    @{Row=29; B="TSTAUTO61"; C="TSTAUTO25"},
    @{Row=30; B="TSTAUTO57"; C="TSTAUTO21"},
    @{Row=31; B="TSTAUTO41"; C="TSTAUTO5"},
    @{Row=32; B="TSTAUTO58"; C="TSTAUTO22"},
    @{Row=33; B="TSTAUTO58"; C="TSTAUTO22"}
)

foreach ($row in $jiraData) {
    $r = $row.Row
    $rng = $jira.Range("B" + $r + ":C" + $r)
    $rng.Font.Color = 0
    $jira.Cells.Item($r, 2).Value = $row.B
    $jira.Cells.Item($r, 3).Value = $row.C
}

# Rows 30-33 now describe the user-management actions (shifted up once the
# duplicate EmployeeCreation/EmployeeUserCreation rows were folded away).
$jira.Range("A30").Value = "ChangeEmployee"
$jira.Range("A31").Value = "CreateUser"
$jira.Range("A32").Value = "ChangeUser"
$jira.Range("A33").Value = "BlockUser"

# The two trailing rows (old ChangeUser/BlockUser) are no longer needed.
$jira.Rows("34:35").Delete()

# Restore the view to the scrolled/selected state captured in the workbook.
$jira.Range("A25").Select()

# ---------------------------------------------------------------------------
# 3) ServerDetails sheet (sheet7) becomes the active tab.
# ---------------------------------------------------------------------------
$server = $wb.Worksheets.Item("ServerDetails")
$server.Activate()
$server.Range("C1").Select()
